# This script applies the updated simulation results for the
# "case with 380 kV done" commit to the pl_mw (line active power) results
# sheet. Columns B, C, E, F, G, L, M, N for rows 2-25 receive new computed
# values; all other cells (A, D, H-K, O, and the header row) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "B2" = 3.285299006035018
    "C2" = 0.3063910664355944
    "E2" = 0.03260438841515523
    "F2" = 0.4443680307746263
    "G2" = 0.002593447952531551
    "L2" = 0.1959427902342341
    "M2" = 0.5588372630720713
    "N2" = 2.913804197205835
    "B3" = 3.153995099164376
    "C3" = 0.2722907654708138
    "E3" = 0.03169665204378269
    "F3" = 0.387822817061874
    "G3" = 0.002600312965820207
    "L3" = 0.1937345121683052
    "M3" = 0.5406136102851136
    "N3" = 2.911298592328734
    "B4" = 3.075547013271944
    "C4" = 0.2514964955564665
    "E4" = 0.03113065065154075
    "F4" = 0.3531389305168915
    "G4" = 0.002604743995309037
    "L4" = 0.1924820619379943
    "M4" = 0.5297840588041751
    "N4" = 2.910435841117931
    "B5" = 3.044120702954274
    "C5" = 0.2430574184021168
    "E5" = 0.03089778326186909
    "F5" = 0.3390132514313251
    "G5" = 0.002606604169846119
    "L5" = 0.1919976269384023
    "M5" = 0.5254608025056626
    "N5" = 2.910252902129017
    "B6" = 3.038934988787275
    "C6" = 0.2416581771890947
    "E6" = 0.03085898065418746
    "F6" = 0.336668177824194
    "G6" = 0.002606916347678512
    "L6" = 0.1919187523306221
    "M6" = 0.5247483396187533
    "N6" = 2.910232673643719
    "B7" = 3.075120998721502
    "C7" = 0.2513825442081554
    "E7" = 0.03112751915422152
    "F7" = 0.3529483938368969
    "G7" = 0.002604768861448339
    "L7" = 0.1924754236932671
    "M7" = 0.5297253907592037
    "N7" = 2.910432692791602
    "B8" = 3.239571617642696
    "C8" = 0.2946027020927318
    "E8" = 0.03229316969212981
    "F8" = 0.4248636149813336
    "G8" = 0.002595770327968627
    "L8" = 0.1951598671836905
    "M8" = 0.5524786827059955
    "N8" = 2.912799251453777
    "B9" = 3.579527885421612
    "C9" = 0.3805596185739546
    "E9" = 0.03451240754549723
    "F9" = 0.5661985755042025
    "G9" = 0.002579827653824767
    "L9" = 0.2012482548512367
    "M9" = 0.5999841097111016
    "N9" = 2.922860582528685
    "B10" = 3.840292000744
    "C10" = 0.444540743035077
    "E10" = 0.03610520015231167
    "F10" = 0.6702781546542269
    "G10" = 0.002569139650201556
    "L10" = 0.2062295139034234
    "M10" = 0.6366939879899363
    "N10" = 2.933641649416131
    "B11" = 3.961387470227521
    "C11" = 0.4738487071786608
    "E11" = 0.03682226621580398
    "F11" = 0.7176906081379002
    "G11" = 0.002564497116735704
    "L11" = 0.2086072741385721
    "M11" = 0.6537976989548184
    "N11" = 2.939300769776509
    "B12" = 4.007604699846127
    "C12" = 0.4849776530816712
    "E12" = 0.03709277244762887
    "F12" = 0.7356546913071611
    "G12" = 0.002562770451888009
    "L12" = 0.2095238380302504
    "M12" = 0.6603333767487101
    "N12" = 2.941553773715185
    "B13" = 3.997634848558732
    "C13" = 0.4825794475465841
    "E13" = 0.0370345591802046
    "F13" = 0.7317853510981394
    "G13" = 0.002563140928244491
    "L13" = 0.2093257195643474
    "M13" = 0.6589231710319581
    "N13" = 2.941063632047502
    "B14" = 3.96518252277923
    "C14" = 0.4747636690634636
    "E14" = 0.03684454131460235
    "F14" = 0.7191683204515869
    "G14" = 0.002564354435668275
    "L14" = 0.2086823560691613
    "M14" = 0.6543342088218367
    "N14" = 2.9394839122684
    "B15" = 3.945351713559603
    "C15" = 0.4699803169339702
    "E15" = 0.03672801703598694
    "F15" = 0.7114413442032514
    "G15" = 0.002565101822184169
    "L15" = 0.2082903841883734
    "M15" = 0.651531026630721
    "N15" = 2.938530660280094
    "B16" = 3.832428344365553
    "C16" = 0.4426296051849476
    "E16" = 0.03605819128505949
    "F16" = 0.6671810134426437
    "G16" = 0.002569447450221141
    "L16" = 0.2060763779219883
    "M16" = 0.6355844175087171
    "N16" = 2.933287126677044
    "B17" = 3.763790635560724
    "C17" = 0.4259037987062584
    "E17" = 0.03564538690424435
    "F17" = 0.6400460337215605
    "G17" = 0.002572169426353806
    "L17" = 0.2047468419851128
    "M17" = 0.6259057428761352
    "N17" = 2.930264749979415
    "B18" = 3.724544436460235
    "C18" = 0.41630253712799
    "E18" = 0.03540724358969882
    "F18" = 0.6244449056556647
    "G18" = 0.002573755706710824
    "L18" = 0.2039926421275311
    "M18" = 0.6203768247513537
    "N18" = 2.928597310310806
    "B19" = 3.711296097978561
    "C19" = 0.4130549289790224
    "E19" = 0.03532648911326497
    "F19" = 0.6191636801734006
    "G19" = 0.002574296350819686
    "L19" = 0.2037390862252266
    "M19" = 0.6185113319949664
    "N19" = 2.928044884961963
    "B20" = 3.771073145709124
    "C20" = 0.4276823141044019
    "E20" = 0.0356894036858133
    "F20" = 0.6429339538360921
    "G20" = 0.002571877529188945
    "L20" = 0.2048872845732461
    "M20" = 0.6269321161370698
    "N20" = 2.930579133249779
    "B21" = 3.974704715879568
    "C21" = 0.4770585084250456
    "E21" = 0.03690038176720734
    "F21" = 0.7228739723492197
    "G21" = 0.002563997149622991
    "L21" = 0.2088708883009502
    "M21" = 0.6556804943765684
    "N21" = 2.939944916466089
    "B22" = 4.109896548103279
    "C22" = 0.5095081176793883
    "E22" = 0.03768583853871554
    "F22" = 0.7751780083420101
    "G22" = 0.002559029571852951
    "L22" = 0.2115686255850591
    "M22" = 0.6748127438659424
    "N22" = 2.946707969919288
    "B23" = 4.037547511916273
    "C23" = 0.4921722250734888
    "E23" = 0.03726715731291463
    "F23" = 0.7472568307916134
    "G23" = 0.002561664210312272
    "L23" = 0.2101201416234204
    "M23" = 0.6645698171748222
    "N23" = 2.943039164327274
    "B24" = 3.767780055055653
    "C24" = 0.4268782018462502
    "E24" = 0.03566950624582077
    "F24" = 0.6416283278902171
    "G24" = 0.002572009429492317
    "L24" = 0.204823758814328
    "M24" = 0.6264679823272701
    "N24" = 2.930436782137207
    "B25" = 3.485654925016263
    "C25" = 0.3571674913833931
    "E25" = 0.03391887560045959
    "F25" = 0.5279251897347308
    "G25" = 0.002583959589928095
    "L25" = 0.199512373534148
    "M25" = 0.5868187296525065
    "N25" = 2.919549755896369
}

foreach ($cellAddress in $newValues.Keys) {
    $ws.Range($cellAddress).Value = $newValues[$cellAddress]
}

Write-Output "Updated $($newValues.Count) cells on sheet '$($ws.Name)'"
